$d = $word.ActiveDocument
$find = $d.Content.Find

# Merge the split "<id>p122v_N</id>" runs back into a single run for each
# of the three occurrences in the document (the fig_p122v_* ids are left
# untouched, matching the diff).
$find.Execute("<id>p122v_1</id>", $false, $false, $false, $false, $false, $true, 1, $false, "<id>p122v_1</id>", 2)
$find.Execute("<id>p122v_2</id>", $false, $false, $false, $false, $false, $true, 1, $false, "<id>p122v_2</id>", 2)
$find.Execute("<id>p122v_3</id>", $false, $false, $false, $false, $false, $true, 1, $false, "<id>p122v_3</id>", 2)
